$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title placeholder: "Testing" + " " + "custom" + " " + "properties" -> one run.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleChars = $titleRange.Characters(1, $titleRange.Length)
$titleChars.Text = "Testing custom properties"

# Subtitle placeholder: two leading line breaks stay untouched; "A." + " " + "M." -> one run.
$subRange = $s.Shapes.Item(2).TextFrame.TextRange
$subChars = $subRange.Characters(3, $subRange.Length - 2)
$subChars.Text = "A. M."
